# "debut partie 19:30 23/12/2014"
# Kills updated for Chine (column G) on the current day's entry row (G39):
# the player typed the new kill count (5, was 4) into G39, which overwrote
# the cell's carry-forward formula (=G40) with a literal value, and the
# selection then moved on to the next cell (G40) as it does after typing
# a value into a cell and confirming it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G39").Select()
$ws.Range("G39").Value = 5
$ws.Range("G40").Select()
